$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.152.03"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "1.636.54"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.65"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("E6").Value = "  +2.05%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.01"
$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").Value = "1.865.74"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").Value = "1.653.11"
$ws.Range("E13").Value = "  +0.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("E14").Value = "  +0.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.541"
$ws.Range("E15").Value = "  +1.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.58"
$ws.Range("E16").Value = "  -0.63%  "

$ws.Range("D17").Value = "27.156.31"
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("E18").Value = "  +1.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.77"
$ws.Range("E19").Value = "  -0.67%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"
$ws.Range("E21").Value = "  +1.35%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.52"
$ws.Range("E23").Value = "  +3.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.48"
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.39"
$ws.Range("E27").Value = "  +2.07%  "

$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0507"
$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  +1.48%  "

$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("D34").Value = "1.300.96"
$ws.Range("E34").Value = "  +3.06%  "

$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +0.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0177"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.545"
$ws.Range("E38").Value = "  +2.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.853"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("E42").Value = "  +5.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.32"
$ws.Range("E43").Value = "  -1.02%  "

$ws.Range("D44").Value = "1.775.60"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.69"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.34"
$ws.Range("E46").Value = "  -0.96%  "

$ws.Range("E47").Value = "  -0.90%  "

$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +1.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"

$ws.Range("E50").Value = "  -0.35%  "

$ws.Range("E51").Value = "  -0.30%  "
